$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 13.458797
$ws.Cells.Item(2, 8).Value = 40.376391
$ws.Cells.Item(2, 9).Value = 0.06830096976102129
$ws.Cells.Item(2, 10).Value = 0.06973720484213804
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 1641.014146666983
$ws.Cells.Item(2, 18).Value = 14769.12732000284
$ws.Cells.Item(2, 19).Value = 0.01558786691678434
$ws.Cells.Item(2, 20).Value = 0.01687445393796308
$ws.Cells.Item(3, 7).Value = 13.458797
$ws.Cells.Item(3, 8).Value = 40.376391
$ws.Cells.Item(3, 9).Value = 0.06830096976102129
$ws.Cells.Item(3, 10).Value = 0.06973720484213804
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 1990.748940861009
$ws.Cells.Item(3, 18).Value = 17916.74046774909
$ws.Cells.Item(3, 19).Value = 0.0189099707750223
$ws.Cells.Item(3, 20).Value = 0.02047075667984781
$ws.Cells.Item(4, 7).Value = 13.458797
$ws.Cells.Item(4, 8).Value = 40.376391
$ws.Cells.Item(4, 9).Value = 0.06830096976102129
$ws.Cells.Item(4, 10).Value = 0.06973720484213804
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 1123.876430748559
$ws.Cells.Item(4, 18).Value = 10114.88787673703
$ws.Cells.Item(4, 19).Value = 0.01067561560575278
$ws.Cells.Item(4, 20).Value = 0.01155675659539424
$ws.Cells.Item(5, 7).Value = 13.458797
$ws.Cells.Item(5, 8).Value = 40.376391
$ws.Cells.Item(5, 9).Value = 0.06830096976102129
$ws.Cells.Item(5, 10).Value = 0.06973720484213804
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 1225.671949889613
$ws.Cells.Item(5, 18).Value = 7354.031699337675
$ws.Cells.Item(5, 19).Value = 0.01164256339734774
$ws.Cells.Item(5, 20).Value = 0.008402342703127976
$ws.Cells.Item(6, 7).Value = 13.458797
$ws.Cells.Item(6, 8).Value = 40.376391
$ws.Cells.Item(6, 9).Value = 0.06830096976102129
$ws.Cells.Item(6, 10).Value = 0.06973720484213804
$ws.Cells.Item(6, 13).Value = 89.83562999999999
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 1209.07950753711
$ws.Cells.Item(6, 18).Value = 10881.71556783399
$ws.Cells.Item(6, 19).Value = 0.01148495306611412
$ws.Cells.Item(6, 20).Value = 0.01243289492580493
$ws.Cells.Item(7, 9).Value = 0.1240039124627887
$ws.Cells.Item(7, 10).Value = 0.1266114708898203
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 2979.345319773504
$ws.Cells.Item(7, 18).Value = 26814.10787796153
$ws.Cells.Item(7, 19).Value = 0.02830057159354778
$ws.Cells.Item(7, 20).Value = 0.03063643629515055
$ws.Cells.Item(8, 9).Value = 0.1240039124627887
$ws.Cells.Item(8, 10).Value = 0.1266114708898203
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 3614.306769898879
$ws.Cells.Item(8, 19).Value = 0.03433202147589382
$ws.Cells.Item(8, 20).Value = 0.03716570830921884
$ws.Cells.Item(9, 9).Value = 0.1240039124627887
$ws.Cells.Item(9, 10).Value = 0.1266114708898203
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 2040.455282335829
$ws.Cells.Item(9, 18).Value = 18364.09754102246
$ws.Cells.Item(9, 19).Value = 0.01938212748214358
$ws.Cells.Item(9, 20).Value = 0.0209818841258513
$ws.Cells.Item(10, 9).Value = 0.1240039124627887
$ws.Cells.Item(10, 10).Value = 0.1266114708898203
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 2225.2702664984
$ws.Cells.Item(10, 18).Value = 13351.6215989904
$ws.Cells.Item(10, 19).Value = 0.02113767077420174
$ws.Cells.Item(10, 20).Value = 0.01525488397436559
$ws.Cells.Item(11, 9).Value = 0.1240039124627887
$ws.Cells.Item(11, 10).Value = 0.1266114708898203
$ws.Cells.Item(11, 13).Value = 89.83562999999999
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 2195.14583669568
$ws.Cells.Item(11, 18).Value = 19756.31253026112
$ws.Cells.Item(11, 19).Value = 0.02085152113700174
$ws.Cells.Item(11, 20).Value = 0.02257255818523405
$ws.Cells.Item(12, 7).Value = 76.51423666666666
$ws.Cells.Item(12, 8).Value = 229.54271
$ws.Cells.Item(12, 9).Value = 0.3882959647030583
$ws.Cells.Item(12, 10).Value = 0.3964610652618627
$ws.Cells.Item(12, 13).Value = 121.928739
$ws.Cells.Item(12, 14).Value = 365.786217
$ws.Cells.Item(12, 15).Value = 0.2282232151508951
$ws.Cells.Item(12, 16).Value = 0.2419720431319445
$ws.Cells.Item(12, 17).Value = 9329.284392314228
$ws.Cells.Item(12, 18).Value = 83963.55953082806
$ws.Cells.Item(12, 19).Value = 0.08861815349465044
$ws.Cells.Item(12, 20).Value = 0.09593249398368012
$ws.Cells.Item(13, 7).Value = 76.51423666666666
$ws.Cells.Item(13, 8).Value = 229.54271
$ws.Cells.Item(13, 9).Value = 0.3882959647030583
$ws.Cells.Item(13, 10).Value = 0.3964610652618627
$ws.Cells.Item(13, 15).Value = 0.2768624053389947
$ws.Cells.Item(13, 16).Value = 0.2935413991166814
$ws.Cells.Item(13, 17).Value = 11317.55205201143
$ws.Cells.Item(13, 18).Value = 101857.9684681029
$ws.Cells.Item(13, 19).Value = 0.1075045547711141
$ws.Cells.Item(13, 20).Value = 0.1163777357922571
$ws.Cells.Item(14, 7).Value = 76.51423666666666
$ws.Cells.Item(14, 8).Value = 229.54271
$ws.Cells.Item(14, 9).Value = 0.3882959647030583
$ws.Cells.Item(14, 10).Value = 0.3964610652618627
$ws.Cells.Item(14, 13).Value = 83.50496933333334
$ws.Cells.Item(14, 14).Value = 250.514908
$ws.Cells.Item(14, 15).Value = 0.1563025480180701
$ws.Cells.Item(14, 16).Value = 0.1657186665504434
$ws.Cells.Item(14, 17).Value = 6389.318986413408
$ws.Cells.Item(14, 18).Value = 57503.87087772068
$ws.Cells.Item(14, 19).Value = 0.06069164866822261
$ws.Cells.Item(14, 20).Value = 0.06570099907436419
$ws.Cells.Item(15, 7).Value = 76.51423666666666
$ws.Cells.Item(15, 8).Value = 229.54271
$ws.Cells.Item(15, 9).Value = 0.3882959647030583
$ws.Cells.Item(15, 10).Value = 0.3964610652618627
$ws.Cells.Item(15, 13).Value = 91.06846250000001
$ws.Cells.Item(15, 14).Value = 182.136925
$ws.Cells.Item(15, 15).Value = 0.1704597085236707
$ws.Cells.Item(15, 16).Value = 0.1204857969594293
$ws.Cells.Item(15, 17).Value = 6968.033892594459
$ws.Cells.Item(15, 18).Value = 41808.20335556675
$ws.Cells.Item(15, 19).Value = 0.06618881696420086
$ws.Cells.Item(15, 20).Value = 0.04776792741145985
$ws.Cells.Item(16, 7).Value = 76.51423666666666
$ws.Cells.Item(16, 8).Value = 229.54271
$ws.Cells.Item(16, 9).Value = 0.3882959647030583
$ws.Cells.Item(16, 10).Value = 0.3964610652618627
$ws.Cells.Item(16, 13).Value = 89.83562999999999
$ws.Cells.Item(16, 14).Value = 269.50689
$ws.Cells.Item(16, 15).Value = 0.1681521229683693
$ws.Cells.Item(16, 16).Value = 0.1782820942415013
$ws.Cells.Item(16, 17).Value = 6873.704654919099
$ws.Cells.Item(16, 18).Value = 61863.3418942719
$ws.Cells.Item(16, 19).Value = 0.06529279080487024
$ws.Cells.Item(16, 20).Value = 0.07068190900010141
$ws.Cells.Item(17, 7).Value = 12.174794
$ws.Cells.Item(17, 8).Value = 24.349588
$ws.Cells.Item(17, 9).Value = 0.06178488588843889
$ws.Cells.Item(17, 10).Value = 0.04205606702633888
$ws.Cells.Item(17, 13).Value = 121.928739
$ws.Cells.Item(17, 14).Value = 365.786217
$ws.Cells.Item(17, 15).Value = 0.2282232151508951
$ws.Cells.Item(17, 16).Value = 0.2419720431319445
$ws.Cells.Item(17, 17).Value = 1484.457280004766
$ws.Cells.Item(17, 18).Value = 8906.743680028596
$ws.Cells.Item(17, 19).Value = 0.01410074530519069
$ws.Cells.Item(17, 20).Value = 0.01017639246445722
$ws.Cells.Item(18, 7).Value = 12.174794
$ws.Cells.Item(18, 8).Value = 24.349588
$ws.Cells.Item(18, 9).Value = 0.06178488588843889
$ws.Cells.Item(18, 10).Value = 0.04205606702633888
$ws.Cells.Item(18, 15).Value = 0.2768624053389947
$ws.Cells.Item(18, 16).Value = 0.2935413991166814
$ws.Cells.Item(18, 17).Value = 1800.82649739802
$ws.Cells.Item(18, 18).Value = 10804.95898438812
$ws.Cells.Item(18, 19).Value = 0.0171059121206685
$ws.Cells.Item(18, 20).Value = 0.01234519675625645
$ws.Cells.Item(19, 7).Value = 12.174794
$ws.Cells.Item(19, 8).Value = 24.349588
$ws.Cells.Item(19, 9).Value = 0.06178488588843889
$ws.Cells.Item(19, 10).Value = 0.04205606702633888
$ws.Cells.Item(19, 13).Value = 83.50496933333334
$ws.Cells.Item(19, 14).Value = 250.514908
$ws.Cells.Item(19, 15).Value = 0.1563025480180701
$ws.Cells.Item(19, 16).Value = 0.1657186665504434
$ws.Cells.Item(19, 17).Value = 1016.655799609651
$ws.Cells.Item(19, 18).Value = 6099.934797657904
$ws.Cells.Item(19, 19).Value = 0.009657135093368699
$ws.Cells.Item(19, 20).Value = 0.00696947534796095
$ws.Cells.Item(20, 7).Value = 12.174794
$ws.Cells.Item(20, 8).Value = 24.349588
$ws.Cells.Item(20, 9).Value = 0.06178488588843889
$ws.Cells.Item(20, 10).Value = 0.04205606702633888
$ws.Cells.Item(20, 13).Value = 91.06846250000001
$ws.Cells.Item(20, 14).Value = 182.136925
$ws.Cells.Item(20, 15).Value = 0.1704597085236707
$ws.Cells.Item(20, 16).Value = 0.1204857969594293
$ws.Cells.Item(20, 17).Value = 1108.739770834225
$ws.Cells.Item(20, 18).Value = 4434.9590833369
$ws.Cells.Item(20, 19).Value = 0.01053183363971155
$ws.Cells.Item(20, 20).Value = 0.005067158752647618
$ws.Cells.Item(21, 7).Value = 12.174794
$ws.Cells.Item(21, 8).Value = 24.349588
$ws.Cells.Item(21, 9).Value = 0.06178488588843889
$ws.Cells.Item(21, 10).Value = 0.04205606702633888
$ws.Cells.Item(21, 13).Value = 89.83562999999999
$ws.Cells.Item(21, 14).Value = 269.50689
$ws.Cells.Item(21, 15).Value = 0.1681521229683693
$ws.Cells.Item(21, 16).Value = 0.1782820942415013
$ws.Cells.Item(21, 17).Value = 1093.73028911022
$ws.Cells.Item(21, 18).Value = 6562.38173466132
$ws.Cells.Item(21, 19).Value = 0.01038925972949944
$ws.Cells.Item(21, 20).Value = 0.007497843705016644
$ws.Cells.Item(22, 7).Value = 70.46836733333333
$ws.Cells.Item(22, 8).Value = 211.405102
$ws.Cells.Item(22, 9).Value = 0.3576142671846927
$ws.Cells.Item(22, 10).Value = 0.36513419197984
$ws.Cells.Item(22, 13).Value = 121.928739
$ws.Cells.Item(22, 14).Value = 365.786217
$ws.Cells.Item(22, 15).Value = 0.2282232151508951
$ws.Cells.Item(22, 16).Value = 0.2419720431319445
$ws.Cells.Item(22, 17).Value = 8592.119168342126
$ws.Cells.Item(22, 18).Value = 77329.07251507913
$ws.Cells.Item(22, 19).Value = 0.08161587784072181
$ws.Cells.Item(22, 20).Value = 0.08835226645069356
$ws.Cells.Item(23, 7).Value = 70.46836733333333
$ws.Cells.Item(23, 8).Value = 211.405102
$ws.Cells.Item(23, 9).Value = 0.3576142671846927
$ws.Cells.Item(23, 10).Value = 0.36513419197984
$ws.Cells.Item(23, 15).Value = 0.2768624053389947
$ws.Cells.Item(23, 16).Value = 0.2935413991166814
$ws.Cells.Item(23, 17).Value = 10423.28134030389
$ws.Cells.Item(23, 18).Value = 93809.53206273497
$ws.Cells.Item(23, 19).Value = 0.09900994619629595
$ws.Cells.Item(23, 20).Value = 0.1071820015791012
$ws.Cells.Item(24, 7).Value = 70.46836733333333
$ws.Cells.Item(24, 8).Value = 211.405102
$ws.Cells.Item(24, 9).Value = 0.3576142671846927
$ws.Cells.Item(24, 10).Value = 0.36513419197984
$ws.Cells.Item(24, 13).Value = 83.50496933333334
$ws.Cells.Item(24, 14).Value = 250.514908
$ws.Cells.Item(24, 15).Value = 0.1563025480180701
$ws.Cells.Item(24, 16).Value = 0.1657186665504434
$ws.Cells.Item(24, 17).Value = 5884.458853140069
$ws.Cells.Item(24, 18).Value = 52960.12967826062
$ws.Cells.Item(24, 19).Value = 0.05589602116858237
$ws.Cells.Item(24, 20).Value = 0.06050955140687268
$ws.Cells.Item(25, 7).Value = 70.46836733333333
$ws.Cells.Item(25, 8).Value = 211.405102
$ws.Cells.Item(25, 9).Value = 0.3576142671846927
$ws.Cells.Item(25, 10).Value = 0.36513419197984
$ws.Cells.Item(25, 13).Value = 91.06846250000001
$ws.Cells.Item(25, 14).Value = 182.136925
$ws.Cells.Item(25, 15).Value = 0.1704597085236707
$ws.Cells.Item(25, 16).Value = 0.1204857969594293
$ws.Cells.Item(25, 17).Value = 6417.445867931892
$ws.Cells.Item(25, 18).Value = 38504.67520759135
$ws.Cells.Item(25, 19).Value = 0.06095882374820883
$ws.Cells.Item(25, 20).Value = 0.0439934841178283
$ws.Cells.Item(26, 7).Value = 70.46836733333333
$ws.Cells.Item(26, 8).Value = 211.405102
$ws.Cells.Item(26, 9).Value = 0.3576142671846927
$ws.Cells.Item(26, 10).Value = 0.36513419197984
$ws.Cells.Item(26, 13).Value = 89.83562999999999
$ws.Cells.Item(26, 14).Value = 269.50689
$ws.Cells.Item(26, 15).Value = 0.1681521229683693
$ws.Cells.Item(26, 16).Value = 0.1782820942415013
$ws.Cells.Item(26, 17).Value = 6330.57017446142
$ws.Cells.Item(26, 18).Value = 6330.570174461420
$ws.Cells.Item(26, 19).Value = 0.06013359823088373
$ws.Cells.Item(26, 20).Value = 0.06509688842534427
